# Update the regression results table in secondary_hypothesis_1 with the
# refreshed analysis numbers (N=251 -> N=265 cohort, re-labeled variable
# names, and updated OR / CI / p-value figures).

$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Row 2: age ---------------------------------------------------------
$t.Cell(2,1).Range.Text  = "Age"
$t.Cell(2,2).Range.Text  = "265"
$t.Cell(2,4).Range.Text  = "0.96, 1.0"
$t.Cell(2,5).Range.Text  = "0.014"

# --- Row 3: sex ----------------------------------------------------------
$t.Cell(3,1).Range.Text  = "Sex"
$t.Cell(3,2).Range.Text  = "265"

# --- Row 5: sex = male ----------------------------------------------------
$t.Cell(5,3).Range.Text  = "1.29"
$t.Cell(5,4).Range.Text  = "0.67, 2.56"

# --- Row 6: bmi ------------------------------------------------------------
$t.Cell(6,1).Range.Text  = "BMI"
$t.Cell(6,2).Range.Text  = "265"
$t.Cell(6,3).Range.Text  = "0.98"
$t.Cell(6,4).Range.Text  = "0.93, 1.02"
$t.Cell(6,5).Range.Text  = "0.3"

# --- Row 7: vis_score --------------------------------------------------
$t.Cell(7,2).Range.Text  = "265"
$t.Cell(7,4).Range.Text  = "1.00, 1.00"
$t.Cell(7,5).Range.Text  = "0.3"

# --- Row 8: pre_cr -----------------------------------------------------
$t.Cell(8,1).Range.Text  = "Creatinine (prior to tMCS)"
$t.Cell(8,2).Range.Text  = "265"
$t.Cell(8,3).Range.Text  = "1.03"
$t.Cell(8,4).Range.Text  = "0.80, 1.30"
$t.Cell(8,5).Range.Text  = "0.8"

# --- Row 9: rrt_group ----------------------------------------------------
$t.Cell(9,1).Range.Text  = "Renal replacement therapy"
$t.Cell(9,2).Range.Text  = "265"
$t.Cell(9,5).Range.Text  = "0.006"

# --- Row 11: RRT before and during tMCS -----------------------------------
$t.Cell(11,3).Range.Text = "0.19"
$t.Cell(11,4).Range.Text = "0.05, 0.55"

# --- Row 12: RRT during tMCS only -----------------------------------------
$t.Cell(12,3).Range.Text = "0.70"
$t.Cell(12,4).Range.Text = "0.35, 1.38"

# --- Row 13: Maximal AKI stadium -> Max KDIGO AKI Stage -------------------
$t.Cell(13,1).Range.Text = "Max KDIGO AKI Stage"
$t.Cell(13,2).Range.Text = "265"

# --- Row 15: S1 ------------------------------------------------------------
$t.Cell(15,3).Range.Text = "0.87"
$t.Cell(15,4).Range.Text = "0.35, 2.07"

# --- Row 16: S2 ------------------------------------------------------------
$t.Cell(16,3).Range.Text = "0.62"
$t.Cell(16,4).Range.Text = "0.22, 1.60"

# --- Row 17: S3 ------------------------------------------------------------
$t.Cell(17,3).Range.Text = "1.11"
$t.Cell(17,4).Range.Text = "0.44, 2.68"

# --- Row 19: model fit statistics footer (merged cell spanning all cols) --
$t.Cell(19,1).Range.Text = "Null deviance = 328; Null df = 264; Log-likelihood = NA; AIC = NA; BIC = NA; Deviance = 309; Residual df = 254; No. Obs. = 265"

Write-Output "done"
